$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Move "Welcome / Sign up" from E6 to G6
$ws.Range("E6").Value = $null
$ws.Range("G6").Value = "Welcome / Sign up"

# Replace "Login" (previously C7) with the new text, now placed at G7
$ws.Range("C7").Value = $null
$ws.Range("G7").Value = "Up & coming Start-Up Companies w/ descriptions  "

# Add the three new rows of TO-DO items
$ws.Range("C9").Value = "Daily news reports"
$ws.Range("C10").Value = "10 days"
$ws.Range("C11").Value = "random array of companies"

# Match the final selection shown in the workbook
$ws.Range("E11").Select() | Out-Null
